$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 holds the "Enterprises density (per 1000 people)" figures for
# Micro (B11), SMEs (C11) and MSMEs (D11). Update the values while keeping
# them stored as text (matching the original shared-string/text cells)
# rather than letting Excel auto-convert them to numbers.
$ws.Range("B11").Value = "'49.23"
$ws.Range("C11").Value = "'4.25"
$ws.Range("D11").Value = "'53.47"

# Setting a value via a leading apostrophe makes Excel apply a "quote
# prefix" style to the cell; restore the original (default) style so the
# cells keep their prior appearance/formatting.
$ws.Range("B11:D11").Style = "Normal"
